$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: environment name, URL, and claim number (litigio related change)
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"
$ws.Range("E3").Value = "'1120170200917   "

# Update the selected cell in the sheet view
$ws.Range("E4").Select() | Out-Null
